$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '69.098.77'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = '3.774.05'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '627.63'
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("E6").Value = '  +1.26%  '
$ws.Range("D7").Value = '3.772.21'
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("D12").Value = '6.77'
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '35.35'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '4.406.49'
$ws.Range("D16").Value = '3.834.99'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '69.116.30'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("E18").Value = '  -3.09%  '
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '7.04'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '467.34'
$ws.Range("D22").Value = '9.55'
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("E23").Value = '  +2.08%  '
$ws.Range("D24").Value = '83.00'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +0.69%  '
$ws.Range("D26").Value = '12.00'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("E27").Value = '  +3.22%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '10.02'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '3.922.07'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").Value = '28.75'
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("E35").Value = '  +19.38%  '
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = '3.723.64'
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").Value = '8.94'
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").Value = '0.966'
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").Value = '154.16'
$ws.Range("D46").Value = '43.26'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = '46.77'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("D51").Value = '1.36'
$ws.Range("E51").Value = '  -1.33%  '
